# Auto-generated edit script applying F-column ("想去人数") updates
# to all four worksheets, per the source diff.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 16
$ws.Range("F8").Value = 1344
$ws.Range("F9").Value = 2031
$ws.Range("F12").Value = 2381
$ws.Range("F13").Value = 639
$ws.Range("F15").Value = 3792
$ws.Range("F18").Value = 2859
$ws.Range("F19").Value = 745
$ws.Range("F22").Value = 66
$ws.Range("F23").Value = 1986
$ws.Range("F24").Value = 1159
$ws.Range("F25").Value = 1758
$ws.Range("F26").Value = 364
$ws.Range("F27").Value = 192
$ws.Range("F28").Value = 7987
$ws.Range("F29").Value = 5478
$ws.Range("F32").Value = 739
$ws.Range("F34").Value = 3459
$ws.Range("F37").Value = 370
$ws.Range("F38").Value = 176
$ws.Range("F40").Value = 4576
$ws.Range("F41").Value = 801
$ws.Range("F42").Value = 50
$ws.Range("F43").Value = 373

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F16").Value = 24
$ws.Range("F19").Value = 63

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 8159
$ws.Range("F3").Value = 351
$ws.Range("F4").Value = 1212

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 8159
$ws.Range("F4").Value = 351
$ws.Range("F5").Value = 1212
$ws.Range("F10").Value = 1344
$ws.Range("F14").Value = 3792
$ws.Range("F16").Value = 2859
$ws.Range("F17").Value = 745
$ws.Range("F20").Value = 1986
$ws.Range("F25").Value = 1159
$ws.Range("F27").Value = 1758
$ws.Range("F29").Value = 364
$ws.Range("F30").Value = 192
$ws.Range("F31").Value = 7987
$ws.Range("F32").Value = 5478
$ws.Range("F33").Value = 63
$ws.Range("F35").Value = 739
$ws.Range("F37").Value = 3459
$ws.Range("F40").Value = 370
$ws.Range("F41").Value = 176
$ws.Range("F44").Value = 4576
$ws.Range("F45").Value = 801
$ws.Range("F46").Value = 50
$ws.Range("F47").Value = 373
